# "Aktualisierung auf nur noch offene Ideen"
# Mark a couple more ideas as "already covered by an existing app" (column C)
# and filter the idea list down to the ones that are still open, i.e. whose
# "Bereits in einer der Apps enthalten" column is blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Hochschulsport" (row 8) and "Stundenplan/Vorlesungsplan..." (row 16, which
# previously only carried a tentative "(x)") are now marked as fully covered.
# This happens before the filter is (re)applied so both rows end up hidden,
# same as the two ideas that already had an "x" in column C.
$ws.Range("C8").Value = "X"
$ws.Range("C16").Value = "X"

# Re-apply the AutoFilter over the idea table, showing only the rows whose
# "Bereits in einer der Apps enthalten" column (C, the 3rd field) is blank -
# i.e. only the ideas that are still open.
$ws.Range("A1:C17").AutoFilter(3, @(""), 7) | Out-Null

# Register/refresh the (hidden) filter-database defined name Excel keeps for
# the active AutoFilter range.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "='HS Hof Android App'!`$A`$1:`$C`$17")
$filterName.Visible = $false

# "App weiterempfehlen" (row 4) gets marked as covered right after, i.e. the
# filter isn't re-run for it, so it keeps showing up despite now also being
# covered.
$ws.Range("C4").Value = "X"

# Leave the selection where the user ended up after filtering the list.
$ws.Range("A17").Select() | Out-Null
